# Scheduled-runner price refresh: update cached market-price / profit
# figures on the Leve profit sheets. Values are plain numeric literals
# (no formulas in this workbook), so each target cell is written directly.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1643.0492
$ws.Range("I15").Value = 1643.0492
$ws.Range("K15").Value = 4929.1476
$ws.Range("M15").Value = -4760.1476
$ws.Range("H19").Value = 1025.2
$ws.Range("I19").Value = 1060.3334
$ws.Range("K19").Value = 1060.3334
$ws.Range("M19").Value = -885.3334
$ws.Range("H64").Value = 3716.5
$ws.Range("H67").Value = 3716.5
$ws.Range("H92").Value = 564.4666999999999
$ws.Range("I92").Value = 550.6667
$ws.Range("K92").Value = 550.6667
$ws.Range("M92").Value = 697.3333
$ws.Range("H96").Value = 453.68182
$ws.Range("I96").Value = 309.1
$ws.Range("K96").Value = 927.3000000000001
$ws.Range("M96").Value = 445.6999999999999
$ws.Range("H100").Value = 40289
$ws.Range("I100").Value = 48536.953
$ws.Range("K100").Value = 48536.953
$ws.Range("M100").Value = -47995.953
$ws.Range("H103").Value = 1286.9231
$ws.Range("I103").Value = 1290
$ws.Range("K103").Value = 3870
$ws.Range("M103").Value = -3284
$ws.Range("H111").Value = 971.3333
$ws.Range("I111").Value = 966.9231
$ws.Range("K111").Value = 2900.7693
$ws.Range("M111").Value = 166.2307000000001
$ws.Range("H116").Value = 8735.049999999999
$ws.Range("I116").Value = 6568.1
$ws.Range("J116").Value = 10902
$ws.Range("K116").Value = 6568.1
$ws.Range("L116").Value = 10902
$ws.Range("M116").Value = -3126.1
$ws.Range("N116").Value = -17786
$ws.Range("H132").Value = 1632.8235
$ws.Range("I132").Value = 1475.037
$ws.Range("J132").Value = 2241.4285
$ws.Range("K132").Value = 4425.111
$ws.Range("L132").Value = 6724.2855
$ws.Range("M132").Value = -1895.111
$ws.Range("N132").Value = -11784.2855
$ws.Range("H137").Value = 12072.538
$ws.Range("I137").Value = 5932.1113
$ws.Range("J137").Value = 17335.762
$ws.Range("K137").Value = 17796.3339
$ws.Range("L137").Value = 52007.28599999999
$ws.Range("M137").Value = -15246.3339
$ws.Range("N137").Value = -57107.28599999999
$ws.Range("H138").Value = 2845.1724
$ws.Range("I138").Value = 2601.8948
$ws.Range("J138").Value = 3307.4
$ws.Range("K138").Value = 7805.6844
$ws.Range("L138").Value = 9922.200000000001
$ws.Range("M138").Value = -2665.6844
$ws.Range("N138").Value = -20202.2

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3655.111
$ws.Range("I74").Value = 2541
$ws.Range("J74").Value = 5883.3335
$ws.Range("K74").Value = 2541
$ws.Range("L74").Value = 5883.3335
$ws.Range("M74").Value = -1667
$ws.Range("N74").Value = -7631.3335
$ws.Range("H77").Value = 3655.111
$ws.Range("I77").Value = 2541
$ws.Range("J77").Value = 5883.3335
$ws.Range("K77").Value = 12705
$ws.Range("L77").Value = 29416.6675
$ws.Range("M77").Value = -8337
$ws.Range("N77").Value = -38152.6675

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1279.1333
$ws.Range("I107").Value = 1312.7858
$ws.Range("K107").Value = 1312.7858
$ws.Range("M107").Value = 607.2141999999999
$ws.Range("H134").Value = 9325.074000000001
$ws.Range("I134").Value = 6581
$ws.Range("J134").Value = 21399
$ws.Range("K134").Value = 19743
$ws.Range("L134").Value = 64197
$ws.Range("M134").Value = -17208
$ws.Range("N134").Value = -69267

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3014.6667
$ws.Range("I31").Value = 2191.05
$ws.Range("J31").Value = 5367.857
$ws.Range("K31").Value = 2191.05
$ws.Range("L31").Value = 5367.857
$ws.Range("M31").Value = -1896.05
$ws.Range("N31").Value = -5957.857
$ws.Range("H34").Value = 3014.6667
$ws.Range("I34").Value = 2191.05
$ws.Range("J34").Value = 5367.857
$ws.Range("K34").Value = 2191.05
$ws.Range("L34").Value = 5367.857
$ws.Range("M34").Value = -1989.05
$ws.Range("N34").Value = -5771.857
$ws.Range("H62").Value = 128680.75
$ws.Range("I62").Value = 146349.42
$ws.Range("K62").Value = 146349.42
$ws.Range("M62").Value = -145725.42
$ws.Range("H65").Value = 128680.75
$ws.Range("I65").Value = 146349.42
$ws.Range("K65").Value = 731747.1000000001
$ws.Range("M65").Value = -728627.1000000001
$ws.Range("H74").Value = 46369.8
$ws.Range("J74").Value = 46369.8
$ws.Range("L74").Value = 46369.8
$ws.Range("N74").Value = -48117.8
$ws.Range("H77").Value = 46369.8
$ws.Range("J77").Value = 46369.8
$ws.Range("L77").Value = 139109.4
$ws.Range("N77").Value = -147845.4
$ws.Range("H107").Value = 643.8570999999999
$ws.Range("I107").Value = 517.4211
$ws.Range("J107").Value = 1845
$ws.Range("K107").Value = 517.4211
$ws.Range("L107").Value = 1845
$ws.Range("M107").Value = 1402.5789
$ws.Range("N107").Value = -5685
$ws.Range("H132").Value = 27284.328
$ws.Range("I132").Value = 18025.686
$ws.Range("K132").Value = 54077.058
$ws.Range("M132").Value = -51547.058
$ws.Range("H134").Value = 49056.652
$ws.Range("I134").Value = 55984.57
$ws.Range("J134").Value = 19959.4
$ws.Range("K134").Value = 167953.71
$ws.Range("L134").Value = 59878.2
$ws.Range("M134").Value = -165418.71
$ws.Range("N134").Value = -64948.2

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5066
$ws.Range("J126").Value = 5316.364
$ws.Range("L126").Value = 15949.092
$ws.Range("N126").Value = -20889.092

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2384.1538
$ws.Range("I61").Value = 2090.4546
$ws.Range("K61").Value = 2090.4546
$ws.Range("M61").Value = -1888.4546
$ws.Range("H104").Value = 13353.8
$ws.Range("J104").Value = 13353.8
$ws.Range("L104").Value = 13353.8
$ws.Range("N104").Value = -20341.8
$ws.Range("H113").Value = 2384.1538
$ws.Range("I113").Value = 2090.4546
$ws.Range("K113").Value = 2090.4546
$ws.Range("M113").Value = 79.54539999999997
$ws.Range("H132").Value = 18524020
$ws.Range("I132").Value = 22227024
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 66681072
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -66678542
$ws.Range("N132").Value = -32060
$ws.Range("H136").Value = 7411516
$ws.Range("I136").Value = 9262921
$ws.Range("K136").Value = 27788763
$ws.Range("M136").Value = -27786213

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 111865.445
$ws.Range("I136").Value = 111865.445
$ws.Range("K136").Value = 335596.335
$ws.Range("M136").Value = -333046.335
